# Auto-generated edit script: update market-price derived columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 175.42857
$ws.Range("I9").Value = 173.8
$ws.Range("J9").Value = 179.5
$ws.Range("K9").Value = 173.8
$ws.Range("L9").Value = 179.5
$ws.Range("M9").Value = -4.800000000000011
$ws.Range("N9").Value = -517.5
# Row 43
$ws.Range("H43").Value = 1949.6666
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 2139.6
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 2139.6
$ws.Range("M43").Value = -931
$ws.Range("N43").Value = -2277.6
# Row 99
$ws.Range("H99").Value = 2326.1177
$ws.Range("I99").Value = 703
$ws.Range("J99").Value = 3211.4546
$ws.Range("K99").Value = 2109
$ws.Range("L99").Value = 9634.363799999999
$ws.Range("M99").Value = -611
$ws.Range("N99").Value = -12630.3638
# Row 107
$ws.Range("H107").Value = 958.7778
$ws.Range("I107").Value = 1363
$ws.Range("J107").Value = 523.46155
$ws.Range("K107").Value = 1363
$ws.Range("L107").Value = 523.46155
$ws.Range("M107").Value = 557
$ws.Range("N107").Value = -4363.46155
# Row 113
$ws.Range("H113").Value = 3164.8
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3164.8
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3164.8
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -9672.799999999999
# Row 132
$ws.Range("H132").Value = 4138.6113
$ws.Range("I132").Value = 3899.7334
$ws.Range("J132").Value = 5333
$ws.Range("K132").Value = 11699.2002
$ws.Range("L132").Value = 15999
$ws.Range("M132").Value = -9169.200199999999
$ws.Range("N132").Value = -21059

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1462
$ws.Range("I2").Value = 1705
$ws.Range("J2").Value = 1300
$ws.Range("K2").Value = 1705
$ws.Range("L2").Value = 1300
$ws.Range("M2").Value = -1592
$ws.Range("N2").Value = -1526
# Row 32
$ws.Range("H32").Value = 1308345.4
$ws.Range("I32").Value = 1497907.6
$ws.Range("J32").Value = 23534.555
$ws.Range("K32").Value = 1497907.6
$ws.Range("L32").Value = 23534.555
$ws.Range("M32").Value = -1497620.6
$ws.Range("N32").Value = -24108.555
# Row 34
$ws.Range("H34").Value = 61448.57
$ws.Range("I34").Value = 55000
$ws.Range("J34").Value = 64028
$ws.Range("K34").Value = 55000
$ws.Range("L34").Value = 64028
$ws.Range("M34").Value = -54729
$ws.Range("N34").Value = -64570
# Row 61
$ws.Range("H61").Value = 8132990
$ws.Range("I61").Value = 12822834
$ws.Range("J61").Value = 3927.4
$ws.Range("K61").Value = 12822834
$ws.Range("L61").Value = 3927.4
$ws.Range("M61").Value = -12822622
$ws.Range("N61").Value = -4351.4
# Row 74
$ws.Range("H74").Value = 12502701
$ws.Range("I74").Value = 1247.4736
$ws.Range("J74").Value = 23813540
$ws.Range("K74").Value = 1247.4736
$ws.Range("L74").Value = 23813540
$ws.Range("M74").Value = -373.4736
$ws.Range("N74").Value = -23815288
# Row 77
$ws.Range("H77").Value = 12502701
$ws.Range("I77").Value = 1247.4736
$ws.Range("J77").Value = 23813540
$ws.Range("K77").Value = 6237.368
$ws.Range("L77").Value = 119067700
$ws.Range("M77").Value = -1869.368
$ws.Range("N77").Value = -119076436
# Row 116
$ws.Range("H116").Value = 1462
$ws.Range("I116").Value = 1705
$ws.Range("J116").Value = 1300
$ws.Range("K116").Value = 1705
$ws.Range("L116").Value = 1300
$ws.Range("M116").Value = 589
$ws.Range("N116").Value = -5888
# Row 122
$ws.Range("H122").Value = 26997.7
$ws.Range("I122").Value = 32003.273
$ws.Range("J122").Value = 3400
$ws.Range("K122").Value = 96009.819
$ws.Range("L122").Value = 10200
$ws.Range("M122").Value = -93559.819
$ws.Range("N122").Value = -15100
# Row 132
$ws.Range("H132").Value = 1542691
$ws.Range("I132").Value = 3876.543
$ws.Range("J132").Value = 5133258
$ws.Range("K132").Value = 11629.629
$ws.Range("L132").Value = 15399774
$ws.Range("M132").Value = -9099.629000000001
$ws.Range("N132").Value = -15404834
# Row 136
$ws.Range("H136").Value = 8132990
$ws.Range("I136").Value = 12822834
$ws.Range("J136").Value = 3927.4
$ws.Range("K136").Value = 38468502
$ws.Range("L136").Value = 11782.2
$ws.Range("M136").Value = -38465952
$ws.Range("N136").Value = -16882.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1462
$ws.Range("I3").Value = 1705
$ws.Range("J3").Value = 1300
$ws.Range("K3").Value = 1705
$ws.Range("L3").Value = 1300
$ws.Range("M3").Value = -1591
$ws.Range("N3").Value = -1528
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -10254
# Row 52
$ws.Range("H52").Value = 57780
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 57780
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 57780
$ws.Range("N52").Value = -58306
# Row 99
$ws.Range("H99").Value = 1654.0869
$ws.Range("I99").Value = 1360.3572
$ws.Range("J99").Value = 2111
$ws.Range("K99").Value = 1360.3572
$ws.Range("L99").Value = 2111
$ws.Range("M99").Value = 137.6428000000001
$ws.Range("N99").Value = -5107
# Row 102
$ws.Range("H102").Value = 39871.285
$ws.Range("I102").Value = 18019.8
$ws.Range("J102").Value = 94500
$ws.Range("K102").Value = 18019.8
$ws.Range("L102").Value = 94500
$ws.Range("M102").Value = -14774.8
$ws.Range("N102").Value = -100990
# Row 121
$ws.Range("H121").Value = 57780
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 57780
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 57780
$ws.Range("N121").Value = -61274
# Row 134
$ws.Range("H134").Value = 3465.4546
$ws.Range("I134").Value = 3159
$ws.Range("J134").Value = 4001.75
$ws.Range("K134").Value = 9477
$ws.Range("L134").Value = 12005.25
$ws.Range("M134").Value = -6942
$ws.Range("N134").Value = -17075.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 439.05884
$ws.Range("I22").Value = 339.33334
$ws.Range("J22").Value = 678.4
$ws.Range("K22").Value = 339.33334
$ws.Range("L22").Value = 678.4
$ws.Range("M22").Value = 10.66665999999998
$ws.Range("N22").Value = -1378.4
# Row 32
$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -684
$ws.Range("N32").ClearContents()
# Row 99
$ws.Range("H99").Value = 2028.4688
$ws.Range("I99").Value = 1801.5
$ws.Range("J99").Value = 2104.125
$ws.Range("K99").Value = 1801.5
$ws.Range("L99").Value = 2104.125
$ws.Range("M99").Value = -303.5
$ws.Range("N99").Value = -5100.125
# Row 105
$ws.Range("H105").Value = 1158.4286
$ws.Range("I105").Value = 1016.3333
$ws.Range("J105").Value = 2011
$ws.Range("K105").Value = 1016.3333
$ws.Range("L105").Value = 2011
$ws.Range("M105").Value = 730.6667
$ws.Range("N105").Value = -5505
# Row 107
$ws.Range("H107").Value = 4464900.5
$ws.Range("I107").Value = 8929009
$ws.Range("J107").Value = 792.2857
$ws.Range("K107").Value = 8929009
$ws.Range("L107").Value = 792.2857
$ws.Range("M107").Value = -8927089
$ws.Range("N107").Value = -4632.2857
# Row 122
$ws.Range("H122").Value = 2007.4445
$ws.Range("I122").Value = 2078
$ws.Range("J122").Value = 1993.3334
$ws.Range("K122").Value = 6234
$ws.Range("L122").Value = 5980.0002
$ws.Range("M122").Value = -3784
$ws.Range("N122").Value = -10880.0002
# Row 126
$ws.Range("H126").Value = 2028.4688
$ws.Range("I126").Value = 1801.5
$ws.Range("J126").Value = 2104.125
$ws.Range("K126").Value = 5404.5
$ws.Range("L126").Value = 6312.375
$ws.Range("M126").Value = -2934.5
$ws.Range("N126").Value = -11252.375

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 19
$ws.Range("H19").Value = 4941
$ws.Range("I19").Value = 1000
# Row 22
$ws.Range("H22").Value = 1001.6
$ws.Range("I22").Value = 650
$ws.Range("J22").Value = 1032.174
$ws.Range("K22").Value = 1950
$ws.Range("L22").Value = 3096.522
$ws.Range("M22").Value = -1781
$ws.Range("N22").Value = -3434.522
# Row 27
$ws.Range("H27").Value = 1001.6
$ws.Range("I27").Value = 650
$ws.Range("J27").Value = 1032.174
$ws.Range("K27").Value = 1950
$ws.Range("L27").Value = 3096.522
$ws.Range("M27").Value = -1848
$ws.Range("N27").Value = -3300.522
# Row 81
$ws.Range("H81").Value = 6957.7334
$ws.Range("I81").Value = 1132.3334
$ws.Range("J81").Value = 8414.083000000001
$ws.Range("K81").Value = 3397.0002
$ws.Range("L81").Value = 25242.249
$ws.Range("M81").Value = -2274.0002
$ws.Range("N81").Value = -27488.249
# Row 84
$ws.Range("H84").Value = 6957.7334
$ws.Range("I84").Value = 1132.3334
$ws.Range("J84").Value = 8414.083000000001
$ws.Range("K84").Value = 10191.0006
$ws.Range("L84").Value = 75726.747
$ws.Range("M84").Value = -4575.000599999999
$ws.Range("N84").Value = -86958.747
# Row 109
$ws.Range("H109").Value = 2264.2856
$ws.Range("I109").Value = 700
$ws.Range("J109").Value = 2890
$ws.Range("K109").Value = 2100
$ws.Range("L109").Value = 8670
$ws.Range("M109").Value = -1060
$ws.Range("N109").Value = -10750
# Row 131
$ws.Range("H131").Value = 6000.1924
$ws.Range("I131").Value = 380
$ws.Range("J131").Value = 7686.25
$ws.Range("K131").Value = 1140
$ws.Range("L131").Value = 23058.75
$ws.Range("M131").Value = 3900
$ws.Range("N131").Value = -33138.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 7923
$ws.Range("I70").Value = 8258.968000000001
$ws.Range("J70").Value = 5840
$ws.Range("K70").Value = 8258.968000000001
$ws.Range("L70").Value = 5840
$ws.Range("M70").Value = -7988.968000000001
$ws.Range("N70").Value = -6380
# Row 73
$ws.Range("H73").Value = 7923
$ws.Range("I73").Value = 8258.968000000001
$ws.Range("J73").Value = 5840
$ws.Range("K73").Value = 8258.968000000001
$ws.Range("L73").Value = 5840
$ws.Range("M73").Value = -7322.968000000001
$ws.Range("N73").Value = -7712
# Row 102
$ws.Range("H102").Value = 1411
$ws.Range("I102").Value = 1393.0834
$ws.Range("J102").Value = 1437.875
$ws.Range("K102").Value = 1393.0834
$ws.Range("L102").Value = 1437.875
$ws.Range("M102").Value = 228.9166
$ws.Range("N102").Value = -4681.875
# Row 122
$ws.Range("H122").Value = 2098.276
$ws.Range("I122").Value = 1884.05
$ws.Range("J122").Value = 2574.3333
$ws.Range("K122").Value = 5652.15
$ws.Range("L122").Value = 7722.999899999999
$ws.Range("M122").Value = -3202.15
$ws.Range("N122").Value = -12622.9999
# Row 126
$ws.Range("H126").Value = 2609
$ws.Range("I126").Value = 2609
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7827
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5357
$ws.Range("N126").ClearContents()
# Row 132
$ws.Range("H132").Value = 3428.348
$ws.Range("I132").Value = 2401.625
$ws.Range("J132").Value = 5775.143
$ws.Range("K132").Value = 7204.875
$ws.Range("L132").Value = 17325.429
$ws.Range("M132").Value = -4674.875
$ws.Range("N132").Value = -22385.429

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 26
$ws.Range("H26").Value = 5145
$ws.Range("I26").Value = 5145
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 5145
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -4850
# Row 34
$ws.Range("H34").Value = 6000
$ws.Range("I34").Value = 6000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 6000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -5828
$ws.Range("N34").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 3050.6667
$ws.Range("I122").Value = 3140.8
$ws.Range("J122").Value = 2600
$ws.Range("K122").Value = 9422.400000000001
$ws.Range("L122").Value = 7800
$ws.Range("M122").Value = -6972.400000000001
$ws.Range("N122").Value = -12700
